# edit.ps1 -- apply the "Dreams" -> "Biology" content rewrite described by the diff.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title -----------------------------------------------------------
Replace-Text "Exploring the Enigmatic Realm of Dreams" "The Marvelous World of Biology: Exploring the Symphony of Life"

# --- Author name: "Emma Williams" -> "Dr. Emily Carter" --------------
# Split into three runs (Dr / . / " Emily Carter") to mirror the source edit.
$authorPara = $d.Paragraphs.Item(2)
$authorStart = $authorPara.Range.Start
$oldLen = "Emma Williams".Length
$authorRange = $d.Range($authorStart, $authorStart + $oldLen)
$authorRange.Text = "Dr"
$afterDr = $d.Range($authorStart + 2, $authorStart + 2)
$afterDr.InsertAfter(".")
$afterDot = $d.Range($authorStart + 3, $authorStart + 3)
$afterDot.InsertAfter(" Emily Carter")

# --- Contact line ------------------------------------------------------
# (Do the longer / more specific replacements first so that "edu" isn't
# matched inside the freshly-inserted "edumail" text.)
Replace-Text "williams@berkeley" "emily725@edumail"
Replace-Text "emma" "carter"
$emailPara = $d.Paragraphs.Item(3)
$emailEnd = $emailPara.Range.End
$tldLen = "edu".Length
# Range.End sits just past the paragraph mark, so back up one character.
$tldRange = $d.Range($emailEnd - 1 - $tldLen, $emailEnd - 1)
$tldRange.Text = "org"

# --- Intro paragraph (font size 24) ------------------------------------
Replace-Text "As humans, we spend a significant portion of our lives in the ethereal realm of dreams, navigating landscapes both familiar and fantastical" "Biology, the study of life, unveils a captivating world teeming with intricate processes, marvelous adaptations, and awe-inspiring diversity"

Replace-Text " Dreams have captivated the imaginations of philosophers, artists, and scientists alike throughout history, inspiring countless works of art, literature, and scientific inquiry" " This scientific discipline invites us to delve into the secrets held within living organisms, from the smallest bacteria to the colossal blue whale"

Replace-Text " From Sigmund Freud's psychoanalytic theories to the latest findings in neuroscience, the quest to understand the nature and significance of dreams continues to be an enthralling intellectual pursuit" " In this journey, we uncover the symphony of life, the delicate balance of ecosystems, and the extraordinary resilience of life on Earth"

Replace-Text "Delving into the labyrinthine corridors of the dreaming mind, we encounter a world governed by its own unique logic, where time and space are fluid, and the boundaries between reality and imagination blur" "Biology unravels the genetic code, the blueprint that guides the development, structure, and function of all living beings"

Replace-Text " Dreams transport us to strange and wonderful places, introduce us to unforgettable characters, and often leave us with lingering emotions and insights upon waking" " We marvel at the intricate dance of molecules, the conversion of energy, and the remarkable symphony of chemical reactions that orchestrate the processes of life"

Replace-Text " These nocturnal journeys have long been a source of fascination, a window into the hidden recesses of our subconscious minds" " Through microscopy, we peer into the hidden world of cells, uncovering their diverse structures and specialized functions, revealing a universe within each minuscule entity"

Replace-Text "Dreams have been the subject of scientific scrutiny for centuries, with researchers employing a variety of methods to unravel their mysteries" "Furthermore, biology illuminates the interdependence of organisms, the delicate web of life that connects all ecosystems"

Replace-Text " From detailed dream diaries to advanced neuroimaging techniques, scientists are gradually piecing together the intricate puzzle of why we dream" " We discover the intricate interactions between species, the delicate balance of predators and prey, and the vital role of biodiversity in maintaining the health of our planet"

Replace-Text " While the precise functions of dreams remain elusive, research suggests that they play a role in memory consolidation, emotional regulation, and creative problem-solving" " As we explore the diversity of life, we gain a deeper understanding of our own place in the intricate tapestry of life on Earth"

# --- Summary heading is unchanged --------------------------------------

# --- Summary paragraph ---------------------------------------------------
Replace-Text "Dreams, those enigmatic nocturnal experiences, have intrigued humanity for millennia" "Biology, the study of life, captivates us with its exploration of the intricate symphony of living organisms, the genetic code that guides their existence, the diversity of ecosystems, and the remarkable resilience of life"

Replace-Text " From artistic and literary musings to scientific investigations, the study of dreams has shed light on the hidden workings of our minds" " It unveils the delicate dance of molecules, the intricate structures of cells, and the interdependence of organisms"

# The last two sentences (plus the period run that separated them) collapse
# into a single trailing sentence before the paragraph's final period.
Replace-Text " While the precise purpose of dreams remains a subject of ongoing research, their role in memory consolidation, emotional regulation, and creative thinking is increasingly recognized. As we continue to probe the depths of the dreaming mind, we may one day come closer to understanding the profound significance of these fleeting yet profound journeys into the realm of the unconscious" " Biology reveals the wonders of adaptation, the marvels of diversity, and the extraordinary resilience of life, expanding our understanding of the interconnectedness of all living beings and inspiring us to appreciate the beauty and fragility of our planet"

# --- Append a new empty paragraph at the very end of the document -------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
